$wb = $excel.ActiveWorkbook

# Report data per language sheet: handoff file base names, xlf names and the
# new "Latest Handback DateTime" to stamp once the handback is recorded.
$reports = @(
    @{
        Sheet       = "zh-cn"
        HandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/591059563fbf49c81c35ce567b6f3532d3ed6d6d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang"
        Rows        = @(
            @{ Row = 2; Md = "53c9dd87-9eb1-4248-aecd-c73be98fef3c.md"; Xlf = "53c9dd87-9eb1-4248-aecd-c73be98fef3c.bacd708eae0fd5a27b6d0a9273ec2e18c6b2cc12.zh-cn.xlf"; HandbackDateTime = "2016-01-08 14:27:33" }
            @{ Row = 3; Md = "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.md"; Xlf = "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.020c6830f98f5756493e16d2d8a9e895c4be6e0f.zh-cn.xlf"; HandbackDateTime = "2016-01-08 14:27:33" }
        )
    },
    @{
        Sheet       = "de-de"
        HandoffBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fc577f1f042e2ba2f6215fa151c7807d5aa9690f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang"
        Rows        = @(
            @{ Row = 2; Md = "53c9dd87-9eb1-4248-aecd-c73be98fef3c.md"; Xlf = "53c9dd87-9eb1-4248-aecd-c73be98fef3c.bacd708eae0fd5a27b6d0a9273ec2e18c6b2cc12.de-de.xlf"; HandbackDateTime = "2016-01-08 14:27:56" }
            @{ Row = 3; Md = "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.md"; Xlf = "9afcfc7c-5f1b-416f-a7e3-c731c38e661d.020c6830f98f5756493e16d2d8a9e895c4be6e0f.de-de.xlf"; HandbackDateTime = "2016-01-08 14:27:56" }
        )
    }
)

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/9f99c5946e13a7c5cc8af41516e78bebf212e335/e2e"

foreach ($report in $reports) {
    $ws = $wb.Worksheets.Item($report.Sheet)

    foreach ($entry in $report.Rows) {
        $row = $entry.Row

        # Status: file has now been handed back.
        $ws.Cells.Item($row, 2).Value = "Handed back"

        # E = Latest Target File (same source .md as column A), with a
        # hyperlink matching the existing column-A link.
        $ws.Cells.Item($row, 5).Value = $entry.Md
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 5), "$mdBase/$($entry.Md)", "", "", $entry.Md) | Out-Null

        # F = Latest Handback File (the handed-back xlf), with a hyperlink
        # matching the existing column-C link.
        $ws.Cells.Item($row, 6).Value = $entry.Xlf
        $ws.Hyperlinks.Add($ws.Cells.Item($row, 6), "$($report.HandoffBase)/$($entry.Xlf)", "", "", $entry.Xlf) | Out-Null

        # G = Latest Handback DateTime: stamp the handback time.
        $ws.Cells.Item($row, 7).Value = $entry.HandbackDateTime
    }
}

# The Overview sheet mirrors the same "Status" text per file/language, so
# bring it in line with the handback status as well.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(2, 2).Value = "Handed back"
$wsOverview.Cells.Item(2, 3).Value = "Handed back"
$wsOverview.Cells.Item(3, 2).Value = "Handed back"
$wsOverview.Cells.Item(3, 3).Value = "Handed back"
